$d = $word.ActiveDocument

# --- Edit 1: "Admins will be inserted..." paragraph ---
# Add a new run after the existing text: " (since admins cannot register through the application)"
$admins = $d.Paragraphs(8)
$aStart = $admins.Range.Start
$aEnd = $admins.Range.End
$aRange = $d.Range($aStart, $aEnd - 1)

$adminsXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Admins will be inserted into the database through some method other than the application</w:t></w:r><w:r><w:t xml:space="preserve"> (since admins cannot register through the application)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$aRange.InsertXML($adminsXml)

# --- Edit 2: "Project names are unique" paragraph ---
# Add " (since there is no other unique attribute" before the existing _GoBack bookmark,
# and ")" right after it, keeping the bookmark intact around the same (zero-width) spot.

# The bookmark sits exactly at the insertion boundary, which makes naive inserts ambiguous
# (they land on the wrong side of the bookmark). Remove it first, rebuild the whole
# paragraph's run content (including the bookmark back in its original spot) in one go.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$proj = $d.Paragraphs(9)
$pStart = $proj.Range.Start
$pEnd = $proj.Range.End
$pRange = $d.Range($pStart, $pEnd - 1)

$projXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Project names are unique</w:t></w:r><w:r><w:t xml:space="preserve"> (since there is no other unique attribute</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$pRange.InsertXML($projXml)

Write-Output "Done. Final text:"
Write-Output $d.Content.Text
